$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1890
$ws1.Range("F5").Value = 46
$ws1.Range("F6").Value = 91
$ws1.Range("F10").Value = 1120
$ws1.Range("F11").Value = 385
$ws1.Range("F13").Value = 77
$ws1.Range("F19").Value = 1263
$ws1.Range("F21").Value = 170
$ws1.Range("F22").Value = 281
$ws1.Range("F24").Value = 614
$ws1.Range("F25").Value = 1043
$ws1.Range("F27").Value = 1949
$ws1.Range("F28").Value = 2420
$ws1.Range("F29").Value = 1192
$ws1.Range("F31").Value = 126
$ws1.Range("F32").Value = 337
$ws1.Range("F33").Value = 574
$ws1.Range("F34").Value = 737
$ws1.Range("F35").Value = 799
$ws1.Range("F36").Value = 102
$ws1.Range("F39").Value = 222
$ws1.Range("F40").Value = 552
$ws1.Range("F41").Value = 646
$ws1.Range("F42").Value = 287

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 39
$ws2.Range("F8").Value = 170
$ws2.Range("F15").Value = 266

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 855

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 855
$ws4.Range("F3").Value = 1890
$ws4.Range("F4").Value = 46
$ws4.Range("F5").Value = 39
$ws4.Range("F6").Value = 91
$ws4.Range("F13").Value = 170
$ws4.Range("F14").Value = 1120
$ws4.Range("F15").Value = 385
$ws4.Range("F17").Value = 77
$ws4.Range("F23").Value = 1263
$ws4.Range("F25").Value = 170
$ws4.Range("F26").Value = 281
$ws4.Range("F28").Value = 1043
$ws4.Range("F29").Value = 2420
$ws4.Range("F31").Value = 1192
$ws4.Range("F35").Value = 126
$ws4.Range("F36").Value = 337
$ws4.Range("F37").Value = 574
$ws4.Range("F40").Value = 737
$ws4.Range("F41").Value = 799
$ws4.Range("F43").Value = 222
$ws4.Range("F44").Value = 552
$ws4.Range("F45").Value = 646
$ws4.Range("F46").Value = 287
